$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting
# Late/Outstanding(heading)/Outstanding columns one to the right.
$ws.Columns("N").Insert()

# The newly inserted column keeps the same character width as its left
# neighbour (column M = 11 characters) rather than inheriting the old
# "Late" column's auto-fit width.
$ws.Columns("N").ColumnWidth = 10.14

# The user ended their session with the "Repayment schedule" tab active
# and cell K16 selected.
$ws.Activate()
$ws.Range("K16").Select()
